# Update "想去人数" (interest count) figures for the 苏州-漫展信息 workbook.
# These same events appear on both the "展览" sheet and the "全部类型" sheet,
# so each value needs to be updated in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$wsExhibit.Range("F3").Value  = 12005
$wsExhibit.Range("F8").Value  = 11908
$wsExhibit.Range("F9").Value  = 499
$wsExhibit.Range("F12").Value = 580
$wsExhibit.Range("F13").Value = 1791
$wsExhibit.Range("F14").Value = 5894
$wsExhibit.Range("F18").Value = 28

# Sheet "全部类型": row -> new F value (same events, different row numbers)
$wsAll.Range("F5").Value  = 12005
$wsAll.Range("F11").Value = 11908
$wsAll.Range("F12").Value = 499
$wsAll.Range("F15").Value = 580
$wsAll.Range("F16").Value = 1791
$wsAll.Range("F18").Value = 5894
$wsAll.Range("F22").Value = 28
